$d = $word.ActiveDocument

# Replace the text of a single run inside a paragraph while preserving
# any other runs already present in that paragraph (e.g. the leading /
# trailing empty <w:r/> "spacer" runs) and preserving the paragraph's own
# <w:pPr/>. InsertXML on the paragraph's own Range (Start..End, which
# excludes the paragraph mark) only replaces the run content - it does
# not touch pPr and does not disturb sibling runs outside that Range.
# NOTE: InsertXML drops any w:rPr placed on the inserted run, so bold /
# italic formatting (if needed) must be re-applied afterwards.
function Set-ParagraphText($para, $newText) {
    $start = $para.Range.Start
    $end = $para.Range.End
    $r = $d.Range($start, $end)
    $runXml = '<?xml version="1.0"?><w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:t xml:space="preserve">' + $newText + '</w:t></w:r>'
    $r.InsertXML($runXml)
}

# Re-apply direct character formatting to the run holding the given text,
# searching ONLY within the supplied paragraph's own range so that sibling
# (empty) runs / other paragraphs are left untouched.
function Set-ParagraphBold($para, $text) {
    $scoped = $d.Range($para.Range.Start, $para.Range.End)
    $scoped.Find.Execute($text, $true, $true, $false, $false, $false, $true, 1, $false)
    if ($scoped.Find.Found) {
        $scoped.Bold = 1
    }
}

function Set-ParagraphItalic($para, $text) {
    $scoped = $d.Range($para.Range.Start, $para.Range.End)
    $scoped.Find.Execute($text, $true, $true, $false, $false, $false, $true, 1, $false)
    if ($scoped.Find.Found) {
        $scoped.Italic = 1
    }
}

# Returns the FIRST paragraph whose text contains the given substring.
function Find-FirstParagraph($matchSubstring) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($matchSubstring)) {
            return $p
        }
    }
    return $null
}

# Returns the LAST paragraph whose text contains the given substring.
function Find-LastParagraph($matchSubstring) {
    $result = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($matchSubstring)) {
            $result = $p
        }
    }
    return $result
}

$oldTitle = "Play Fortune Tellers Charm Free - Exciting Magic & Fortune Theme"
$newTitle = "Play Fortune Tellers Charm Free: Review of Leander Games Slot"

# Main page title (Heading1) - the first paragraph in the document carrying
# this text.
Set-ParagraphText (Find-FirstParagraph $oldTitle) $newTitle

# What we like bullets
Set-ParagraphText (Find-FirstParagraph "Easy game mechanics and well-designed test mode") "Easy gameplay mechanics and user-friendly interface"
Set-ParagraphText (Find-FirstParagraph "Immeasurable opportunities to win and lucrative bonus modes") "Abundant opportunities to win and exciting bonus modes"
Set-ParagraphText (Find-FirstParagraph "Beautiful symbol design and immersive fortune teller theme") "Beautiful symbol design and immersive theme"
Set-ParagraphText (Find-FirstParagraph "Unique and diverse theme compared to other virtual slot machines") "Unique and diverse from other slot machines"

# What we don't like bullets
Set-ParagraphText (Find-FirstParagraph "Limited number of pay lines compared to other slot games") "Limited number of pay lines"
Set-ParagraphText (Find-FirstParagraph "Absence of progressive jackpot") "No progressive jackpot feature"

# Bold "title" line near the end of the document (repeats the heading
# text, bold-formatted) - the LAST paragraph carrying the old title text.
$boldPara = Find-LastParagraph $oldTitle
Set-ParagraphText $boldPara $newTitle
Set-ParagraphBold (Find-LastParagraph $newTitle) $newTitle

# Italic meta-description line that immediately follows the bold title
# line, at the very end of the document.
$oldMeta = "Try your fortune with Fortune Tellers Charm, a free slot game by Leander Games, featuring a beautiful fortune teller theme and diverse game mechanics."
$newMeta = "Play Fortune Tellers Charm for free and enjoy an immersive gaming experience. Review of Leander Games slot."
$italicPara = Find-LastParagraph $oldMeta
Set-ParagraphText $italicPara $newMeta
Set-ParagraphItalic (Find-LastParagraph $newMeta) $newMeta
